# lock the excel title row, change the attached file attribute
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Upload File")

# --- Change the attached file attribute ---
# "Attached File" currently lives at the end (col Z). Move it to sit right
# after the "Question Is Active" column (col J), i.e. make it column K, by
# inserting a new column at K and deleting the now-redundant trailing column.
$ws.Range("K1").EntireColumn.Insert()
$ws.Range("K1").Value = "Attached File"
$ws.Range("AA1").EntireColumn.Delete()

# --- Lock the title row ---
# Columns A:D (core identifying fields) stay locked, no special alignment.
$ws.Range("A1:D1").Locked = $true

# Columns E:K (remaining single-answer fields incl. the relocated Attached
# File) stay locked and vertically centered.
$ws.Range("E1:K1").Locked = $true
$ws.Range("E1:K1").VerticalAlignment = -4108

# Columns L:Z (the repeated "answer" blocks) remain unlocked/editable but
# keep the vertical centering.
$ws.Range("L1:Z1").Locked = $false
$ws.Range("L1:Z1").VerticalAlignment = -4108

# Protect the worksheet so the locked title row can't be edited.
$ws.Protect("")

# Restore the selection to where the user left off editing.
$ws.Range("I12").Select()
